$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Table on slide 16 ("PLENARY - COMPLETE THE MISSING GAPS"): switch the
#    table's style from the custom "Table_0" style to the built-in style
#    {6911C6A6-B2E2-47C0-92FC-D5883CB8A8BF}.
# ---------------------------------------------------------------------------
$slide16 = $p.Slides.Item(16)
for ($i = 1; $i -le $slide16.Shapes.Count; $i++) {
    $shp = $slide16.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{6911C6A6-B2E2-47C0-92FC-D5883CB8A8BF}")
    }
}

# ---------------------------------------------------------------------------
# 2) Theme: recolour the deck's theme color scheme from the "Integral" set
#    back to the stock "Office Theme" colors (dk1/lt1/dk2/lt2/accent1-6/
#    hlink/folHlink), matching the target clrScheme.
# ---------------------------------------------------------------------------
$master = $p.Slides.Item(1).Master
$colorScheme = $master.Theme.ThemeColorScheme

# Index order matches MsoThemeColorSchemeIndex:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
# 8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
# Values are the classic VBA RGB() encoding (R + G*256 + B*65536).
$officeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

for ($i = 1; $i -le 12; $i++) {
    $colorScheme.Item($i).RGB = $officeColors[$i - 1]
}
